$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content first
$ws.Cells.Clear()

# Column widths (characters)
$ws.Columns.Item(1).ColumnWidth = 13.33
$ws.Columns.Item(2).ColumnWidth = 41.67
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(4).ColumnWidth = 129.0
$ws.Columns.Item(5).ColumnWidth = 24.67
$ws.Columns.Item(7).ColumnWidth = 11.83

# Cell values
$ws.Range('A1').Value = 'Agências'
$ws.Range('B1').Value = 'Movimentações'
$ws.Range('C1').Value = 'Movimentações'
$ws.Range('D1').Value = 'Status'
$ws.Range('A2').Value = 'Incluir'
$ws.Range('B2').Value = 'ok'
$ws.Range('C2').Value = 'Incluir'
$ws.Range('D2').Value = 'Instrução insert conflitou com a restrição do foreign key "FK_MOVIMENTACAO_HIS". O conflito ocorreu na tabela "dbo.HISTORICOS", column ''ID_HIS'' '
$ws.Range('A3').Value = 'Alterar'
$ws.Range('B3').Value = 'ok'
$ws.Range('C3').Value = 'Alterar'
$ws.Range('D3').Value = 'Falha ao converter data e/ou hora da cadeia de caracacteres'
$ws.Range('A4').Value = 'Excluir'
$ws.Range('B4').Value = 'não é possível deletar, tabela relacionada'
$ws.Range('C4').Value = 'Excluir'
$ws.Range('D4').Value = 'ok'
$ws.Range('A6').Value = 'Clientes'
$ws.Range('C6').Value = 'Contas Correntes'
$ws.Range('A7').Value = 'Incluir'
$ws.Range('B7').Value = 'ok'
$ws.Range('C7').Value = 'Incluir'
$ws.Range('D7').Value = 'ok'
$ws.Range('A8').Value = 'Alterar'
$ws.Range('B8').Value = 'ok'
$ws.Range('C8').Value = 'Alterar'
$ws.Range('D8').Value = 'Conflito foreign key "fk_contacorrente_cli". O conflito ocorreu na tabela "dbo.CLIENTES", column ''ID_CLI'''
$ws.Range('A9').Value = 'Excluir'
$ws.Range('B9').Value = 'não é possível deletar, tabela relacionada'
$ws.Range('C9').Value = 'Excluir'
$ws.Range('D9').Value = 'não é possível deletar, tabela relacionada'
$ws.Range('A11').Value = 'Históricos'
$ws.Range('C11').Value = 'Usuários'
$ws.Range('A12').Value = 'Ler'
$ws.Range('B12').Value = 'ok'
$ws.Range('C12').Value = 'Incluir'
$ws.Range('D12').Value = 'insert conflitou com a restrição foreign key "fk_usuarios_cc". O conflito ocorreu na tabela "dbo.CONTACORRENTE"'
$ws.Range('C13').Value = 'Alterar'
$ws.Range('D13').Value = 'ok'
$ws.Range('C14').Value = 'Excluir'
$ws.Range('D14').Value = 'ok'
$ws.Range('A16').Value = 'Funcionários'
$ws.Range('A17').Value = 'Incluir'
$ws.Range('B17').Value = 'ok'
$ws.Range('A18').Value = 'Alterar'
$ws.Range('B18').Value = 'ok'
$ws.Range('A19').Value = 'Excluir'
$ws.Range('B19').Value = 'ok'

# Styles: green (ok) = style index 2 fill FF00B050, red (problem) = style index 1 fill FFFF0000
$greenCells = 'A2','A3','C4','A7','C7','A8','A12','C13','C14','A17','A18','A19'
foreach ($addr in $greenCells) { $ws.Range($addr).Interior.Color = 5287936 }
$redCells = 'C2','C3','A4','C8','A9','C9','C12'
foreach ($addr in $redCells) { $ws.Range($addr).Interior.Color = 255 }

$ws.Range("A8").Select() | Out-Null
